# Generate Report for Handoff
# Updates the "b.md" row on each sheet (Overview, zh-cn, de-de) to reflect
# that the file is now ready for handoff, with a fresh handoff package and
# an error detail noting the handback file is stale.
#
# NOTE: several of the string values below (e.g. "False", timestamps)
# look like booleans/dates to Excel's type inference, so they are written
# with a leading "'" to force plain text and avoid being auto-coerced into
# a Boolean/Date cell type. Excel itself strips the leading apostrophe from
# the stored/read value.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f10cd102aa574befd68537b82d1279a21e242713/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2017de319ad69a5bd56f4f82b9bf3d1805f45fb2/e2e/b.md."

# ---- Overview sheet: row 3 is the b.md summary row ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value() = "Ready for handoff"
$wsOverview.Range("F3").Value() = "Ready for handoff"
$wsOverview.Range("G3").Value() = "2016-08-16 20:34:17"

# ---- zh-cn sheet: row 3 is the b.md detail row ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value() = "Ready for handoff"
$wsZhCn.Range("F3").Value() = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value() = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value() = "2016-08-16 20:34:12"
$wsZhCn.Range("P3").Value() = "'" + $errorDetail
$wsZhCn.Range("P3").Style = "Normal"
$wsZhCn.Columns.Item(16).ColumnWidth = 40

# ---- de-de sheet: row 3 is the b.md detail row ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value() = "Ready for handoff"
$wsDeDe.Range("F3").Value() = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value() = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value() = "2016-08-16 20:34:17"
$wsDeDe.Range("P3").Value() = "'" + $errorDetail
$wsDeDe.Range("P3").Style = "Normal"
$wsDeDe.Columns.Item(16).ColumnWidth = 40
